$d = $word.ActiveDocument

# 1. Remove the "Existing loans" heading paragraph that precedes the FIRST
#    existing-loans table (first applicant's section). This paragraph also
#    carries a <w:sectPr> (page section break) which is removed along with it.
#    (We walk the Paragraphs collection to reliably find the FIRST match,
#    since locating it through a Find-result Range's own .Paragraphs
#    collection is unreliable in this engine.)
$targetIndex = -1
$i = 1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Existing loans*") {
        $targetIndex = $i
        break
    }
    $i = $i + 1
}
if ($targetIndex -gt 0) {
    $d.Paragraphs.Item($targetIndex).Range.Delete()
}

# 2. Remove the empty paragraph (carrying a continuous <w:sectPr>) that sits
#    right after the SECOND existing-loans table (second applicant's section).
#    NOTE: this, and the deletion above, must happen BEFORE any in-place
#    Range.Text rewriting (step 3 below), because rewriting the text of a
#    sub-paragraph Range leaves this engine's Paragraphs collection reporting
#    stale/incorrect text on every subsequent query.
$cnt = $d.Paragraphs.Count
$sectBreakPara = $d.Paragraphs.Item($cnt - 1)
$sectBreakPara.Range.Delete()

# 3. Remove the <w:lastRenderedPageBreak/> stored before "Bank" in the FIRST
#    existing-loans table (Table 2). Re-writing the cell's text forces the
#    engine to regenerate the run without the stale rendering marker.
$tbl1 = $d.Tables.Item(2)
$bankPara = $tbl1.Cell(1, 1).Range.Paragraphs.Item(1)
$bankRange = $d.Range($bankPara.Range.Start, $bankPara.Range.End - 1)
$bankRange.Text = "BankXX"
$bankRange2 = $d.Range($bankPara.Range.Start, $bankPara.Range.Start + 6)
$bankRange2.Text = "Bank"

# 4. Change the loan amounts to spelled-out numbers (verbalize plugin demo).
$d.Content.Find.Execute("10000", $true, $false, $false, $false, $false, $true, 1, $false, "ten thousand", 2) | Out-Null
$d.Content.Find.Execute("2000", $true, $false, $false, $false, $false, $true, 1, $false, "two thousand", 2) | Out-Null
